$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column F (dSF) per row, as described by the diff.
$values = @{
    2  = -1
    3  = -1
    5  = -1
    6  = -3
    7  = -3
    8  = 2
    9  = 2
    10 = -2
    12 = -3
    13 = -1
    14 = 2
    15 = -2
    16 = 5
    17 = 1
    18 = 1
    19 = 0
    20 = 7
    21 = -2
    22 = 7
    23 = 1
    24 = -3
    25 = -3
    26 = -4
    27 = 5
    28 = -2
    29 = 3
    30 = 2
    31 = -3
    32 = 3
    33 = 2
    34 = 0
    35 = 2
    36 = -4
    37 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
